$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Corrected "Saturate" (MaxSaturates, column F) values for every age row
#    (2-31). This is the bulk of the "corrected Saturate values" part of the
#    commit.
# ---------------------------------------------------------------------------
$saturates = @{
    2  = 3.3
    3  = 4
    4  = 5.3
    5  = 6.67
    6  = 7.33
    7  = 8
    8  = 9.33
    9  = 10.67
    10 = 12.67
    11 = 17.33
    12 = 20
    13 = 21.33
    14 = 20
    15 = 19.33
    16 = 19.33
    17 = 18.67
    18 = 18.67
    19 = 18
    20 = 18
    21 = 17.33
    22 = 16.67
    23 = 16.67
    24 = 16
    25 = 15.33
    26 = 15.33
    27 = 14.67
    28 = 14.67
    29 = 14
    30 = 14
    31 = 13.33
}

foreach ($row in $saturates.Keys) {
    $ws.Cells.Item($row, 6).Value = $saturates[$row]
}

# ---------------------------------------------------------------------------
# 2. "Fixed Flawless Mechanic" - re-apply the wrap-text formatting on the
#    merged "Ages" label cells in column B. The age-group blocks keep the
#    same wrapping behaviour they had before (rows 20-22 wrapped, every other
#    block not wrapped) but the underlying style bookkeeping gets refreshed,
#    which is what re-asserting the formatting here reproduces.
# ---------------------------------------------------------------------------
$ws.Range("B2:B19").WrapText = $false
$ws.Range("B20:B22").WrapText = $true
$ws.Range("B23:B31").WrapText = $false

# ---------------------------------------------------------------------------
# 3. Update the remembered selection to match where the author ended up
#    (cell F32, just below the last data row) before saving.
# ---------------------------------------------------------------------------
[void]$ws.Range("F32").Select()
